$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.685.59'
$ws.Range("E2").Value = '  +0.88%  '
$ws.Range("D3").Value = '1.645.08'
$ws.Range("E3").Value = '  +1.12%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '215.68'
$ws.Range("E5").Value = '  +1.21%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.505'
$ws.Range("E6").Value = '  +1.35%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.01'
$ws.Range("E7").Value = '  +0.15%  '
$ws.Range("E8").Value = '  +1.31%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0627'
$ws.Range("E9").Value = '  +0.62%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.13'
$ws.Range("E10").Value = '  +0.91%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0844'
$ws.Range("E11").Value = '  +0.01%  '
$ws.Range("D12").Value = '1.871.70'
$ws.Range("D13").Value = '1.629.52'
$ws.Range("E13").Value = '  +1.01%  '
$ws.Range("E14").Value = '  +0.91%  '
$ws.Range("E15").Value = '  +1.80%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.04'
$ws.Range("E16").Value = '  +0.74%  '
$ws.Range("D17").Value = '26.689.82'
$ws.Range("E17").Value = '  +0.94%  '
$ws.Range("E18").Value = '  +0.49%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '216.89'
$ws.Range("E19").Value = '  +0.75%  '
$ws.Range("E20").Value = '  +0.13%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.36'
$ws.Range("E21").Value = '  +1.19%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.27'
$ws.Range("E22").Value = '  +0.45%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.49'
$ws.Range("E23").Value = '  +1.91%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.26'
$ws.Range("E24").Value = '  +13.31%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '145.53'
$ws.Range("E25").Value = '  -1.54%  '
$ws.Range("E26").Value = '  +0.29%  '
$ws.Range("E27").Value = '  +0.37%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.13'
$ws.Range("E28").Value = '  +4.24%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.73'
$ws.Range("E29").Value = '  +1.08%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0514'
$ws.Range("E30").Value = '  +1.03%  '
$ws.Range("E31").Value = '  +1.64%  '
$ws.Range("E32").Value = '  +1.49%  '
$ws.Range("E33").Value = '  +2.14%  '
$ws.Range("D34").Value = '1.278.48'
$ws.Range("E34").Value = '  +4.93%  '
$ws.Range("E35").Value = '  +3.44%  '
$ws.Range("E36").Value = '  +1.54%  '
$ws.Range("E37").Value = '  +3.17%  '
$ws.Range("E38").Value = '  +6.06%  '
$ws.Range("E39").Value = '  +3.26%  '
$ws.Range("E40").Value = '  +0.19%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.812'
$ws.Range("E41").Value = '  +2.33%  '
$ws.Range("E42").Value = '  -0.09%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.44'
$ws.Range("E43").Value = '  +1.29%  '
$ws.Range("D44").Value = '1.781.85'
$ws.Range("E44").Value = '  +1.11%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '91.79'
$ws.Range("E45").Value = '  -1.10%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '59.55'
$ws.Range("E46").Value = '  +8.71%  '
$ws.Range("E47").Value = '  +1.30%  '
$ws.Range("E48").Value = '  +1.24%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.75'
$ws.Range("E49").Value = '  +3.10%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0966'
$ws.Range("E50").Value = '  +1.54%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.408'
$ws.Range("E51").Value = '  +0.12%  '
